$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 attribution values updated to reflect relative direction (per commit message).
$rowNum = 2
$updates = [ordered]@{
    "A" = -0
    "B" = -0.107670151431414
    "C" = -0
    "D" = 0.2053876115891933
    "E" = -0.000870470479885202
    "G" = 0
    "J" = -0
    "K" = -0.01825145604069086
    "L" = -0
    "M" = 0.2064948814758124
    "N" = 0.00208025141448778
    "R" = -0
    "S" = 0
    "T" = -0.09778338212392844
    "V" = 0.01432430740601758
    "W" = -0.04174952113749071
    "AB" = 0
    "AC" = -0.07844340024108128
    "AD" = 0
    "AE" = -0.01123724989979894
    "AF" = 0.00467995902030148
    "AG" = -0
    "AI" = -0
    "AJ" = 0
    "AK" = -0
    "AL" = -0.03477317437924476
    "AM" = 0
    "AN" = 0.02975591359928845
    "AO" = 0.06894031625633901
    "AR" = -0
    "AT" = 0
    "AU" = -0.1504158949591624
    "AW" = 0.06762809990804736
    "AX" = 0.01316411438484101
    "AY" = -0
    "BC" = -0
    "BD" = -0.03736677967854636
    "BF" = 0.07026386262018187
    "BG" = 0.02677008819260617
    "BI" = -0
    "BJ" = -0
    "BL" = 0
    "BM" = 0.01906614199068762
    "BO" = -0.05327664047297926
    "BP" = -0.06915173979540963
    "BU" = 0
    "BV" = -0.01804737642947697
    "BX" = 0.01512898727912453
    "BY" = -0.005352203034577101
    "BZ" = -0
    "CB" = 0
    "CD" = 0
    "CE" = 0.0336355770062513
    "CG" = -0.03569317361011079
    "CH" = 0.0251507005628732
    "CJ" = -0
    "CM" = -0
    "CN" = -0.01586600892987506
    "CO" = -0
    "CP" = 0.02594067695371425
    "CQ" = 0.05198710927025178
    "CT" = 0
    "CU" = -0
    "CV" = -0
    "CW" = 0.04516854268119422
    "CY" = -0.03815737946181457
    "CZ" = 0.01954564062428951
    "DB" = 0
    "DE" = -0
    "DF" = 0.04358123142453846
    "DH" = 0.01717012150273894
    "DI" = 0.04988333764285281
    "DJ" = 0
    "DK" = -0
    "DL" = -0
    "DN" = 0
    "DO" = -0.02121619411584265
    "DQ" = 0.03359901195050785
    "DR" = -0.01385797544263802
    "DS" = -0
    "DW" = 0
    "DX" = -0.05212903656296591
    "DY" = -0
    "DZ" = -0.004971154324998677
    "EA" = -0.02718932500684638
    "EB" = 0
    "EF" = -0
    "EG" = 0.05015491900709527
    "EI" = 0.07613886383889099
    "EJ" = -0.0199689068547195
    "EO" = 0
    "EP" = 0.04747992163333337
    "EQ" = 0
    "ER" = -0.03421357883596667
    "ES" = 0.04056784176355664
    "ET" = 0
    "EU" = -0
    "EV" = 0
    "EX" = 0
    "EY" = 0.047632734695192
    "FA" = -0.03127982317017986
    "FB" = 0.02441504420193668
    "FD" = -0
    "FF" = -0
    "FG" = -0
    "FH" = -0.0002063141115742818
    "FI" = 0
    "FJ" = -0.01377775761390556
    "FK" = -0.001886070509032845
    "FL" = -0
    "FP" = -0
    "FQ" = -0.01316388935840617
    "FR" = -0
    "FS" = -0.0009316850821103231
    "FT" = 0.0192548302470746
    "FV" = -0
    "FW" = -0
    "FY" = 0
    "FZ" = -0.03444598377532641
    "GA" = -0
    "GB" = 0.04121121055052257
    "GD" = 0
    "GE" = -0
}

foreach ($col in $updates.Keys) {
    $ws.Range("$col$rowNum").Value = $updates[$col]
}